# Iran weekly-deaths prediction workbook update:
# Insert a new week row ("2021-01-02" / "03 Jan -- 09 Jan 2021") before the
# existing last row (which held the "2021-01-09" / "10 Jan -- 16 Jan 2021" week),
# pushing that row down from row 50 to row 51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 50; existing row 50 (and anything below)
# shifts down to row 51.
$ws.Rows.Item(50).EntireRow.Insert()

# Populate the newly inserted row 50 with the new weekly prediction data.
$ws.Range("A50").Value = "2021-01-02"
$ws.Range("B50").Value = "03 Jan -- 09 Jan 2021"
$ws.Range("C50").Value = 94.56999999999999
$ws.Range("D50").Value = 133.81
$ws.Range("E50").Value = 39.24
$ws.Range("F50").Value = "KNN"
$ws.Range("G50").Value = 1.11
$ws.Range("H50").Value = 33.79
$ws.Range("I50").Value = 41.71
$ws.Range("J50").Value = 45.26
$ws.Range("K50").Value = 41.52
